# Fix issue #51 move keep variables after meta variables.
# Add a new summary row (row 3) for date "11/18/16" with its frequency stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the date label as literal text (matching the existing "9/9/16" text
# cell above it) rather than letting Excel auto-convert it to a date serial.
$dateCell = $ws.Cells.Item(3, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "11/18/16"
$dateCell.ClearFormats()

$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(3, 3).Value = 20
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 20
